# new_ph_pl.xlsx — "day" sheet: fix D50:D52 typing + append 4 freshly
# scraped rows (53-56), per "break out stock.yaml completed".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- D50:D52 (bsecode) were written as text; store them as numbers instead ---
$ws.Range("D50").Value = 532540
$ws.Range("D51").Value = 526299
$ws.Range("D52").Value = 524804

# --- Append 4 new rows (53-56) with freshly scraped stock data ---
# sr, nsecode, name, bsecode, per_chg, close, volume, timeframe, Date Time
$newRows = @(
    @(1, "DRREDDY",    "Dr. Reddy's Laboratories Limited", "500124", 0.39,  6078.4,  395506,  "day", "25/06/2024 11:35:29"),
    @(2, "INDIGO",     "Interglobe Aviation Limited",       "539448", -1.9,  4233.5,  1512109, "day", "25/06/2024 11:35:29"),
    @(3, "LALPATHLAB", "Dr. Lal Path Labs Ltd.",            "539524", 1.09,  2718.95, 253428,  "day", "25/06/2024 11:35:29"),
    @(4, "LUPIN",      "Lupin Limited",                     "500257", -0.5,  1558.85, 493371,  "day", "25/06/2024 11:35:29")
)

$row = 53
foreach ($r in $newRows) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]

    # bsecode must stay text (it's zero-padded-style data elsewhere), so
    # force a text number format before the write, then drop the
    # formatting back to the sheet's default so no stray style lingers.
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $r[3]
    $dCell.ClearFormats()

    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = $r[8]
    $row++
}
